$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that currently sits
#    right under the main H1 title.
# ------------------------------------------------------------------
$metaRange = $d.Content
$foundMeta = $metaRange.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundMeta) {
    $metaRange.Expand(4)   # wdParagraph - grab the whole paragraph, incl. the mark
    $metaRange.Delete()
}

# ------------------------------------------------------------------
# 2. Insert a new bold paragraph ("Play East Sea Dragon King for Free
#    - NetEnt Slot Game") right before the closing image-prompt
#    paragraph.
# ------------------------------------------------------------------
$imgRange = $d.Content
$foundImg = $imgRange.Find.Execute("Create a feature image", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$imgRange.Expand(4)
$startOfImgPara = $imgRange.Start

$titleText = "Play East Sea Dragon King for Free – NetEnt Slot Game"

$insertPoint = $d.Range($startOfImgPara, $startOfImgPara)
$insertPoint.InsertBefore($titleText)

$splitPoint = $d.Range($startOfImgPara + $titleText.Length, $startOfImgPara + $titleText.Length)
$splitPoint.InsertParagraphAfter()

$titleRange = $d.Range($startOfImgPara, $startOfImgPara + $titleText.Length)
$titleRange.Font.Bold = 1

# ------------------------------------------------------------------
# 3. Replace the old image-generation prompt text with the new meta
#    description text, keeping the paragraph's italic formatting.
# ------------------------------------------------------------------
$oldText = 'Create a feature image for the game "East Sea Dragon King" that captures the Asian and cartoon theme of the game. The image should prominently feature a happy Maya warrior with glasses, fitting in with the overall aesthetic of the game. The image should be bright and eye-catching, with a colour palette that reflects the underwater setting and elements of traditional Asian design. The overall style should be cartoonish and fun, with a strong emphasis on the character of the Maya warrior.'
$newText = 'Discover East Sea Dragon King, the newest NetEnt Slot Game. Play for free with high volatility, excellent payouts, and dynamic features. No downloads needed.'

$replaceRange = $d.Content
$replaceRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
